$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = (Get-Date -Year 2016 -Month 8 -Day 25 -Hour 21 -Minute 20 -Second 30)
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

$ws.Range("B6").Value = -18
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 43
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 92
$ws.Range("G6").Value = 42346
$ws.Range("H6").Value = 22079
$ws.Range("I6").Value = 1329
$ws.Range("J6").Value = 194
$ws.Range("K6").Value = 150
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = 24
$ws.Range("N6").Value = "Named"
